$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Fecha de firma: el texto " _____________________" + "_  " + "/" +
#    "202" (envuelto en marcas <w:proofErr> de gramatica) se limpia y
#    se reescribe como " ______________________  " + "/202" sin los
#    artefactos de revision gramatical, dejando la fecha lista para
#    imprimir.
# -----------------------------------------------------------------
$anchor = $d.Content
$foundAnchor = $anchor.Find.Execute("En constancia de la veracidad de la informaci", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundAnchor) { throw "No se encontro el parrafo de la fecha de firma." }

$dateZone = $d.Range($anchor.End, $d.Content.End)
$foundDate = $dateZone.Find.Execute(" _____________________", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundDate) { throw "No se encontro el bloque de guiones de la fecha." }

# El bloque viejo cubre " _____________________" + "_  " + "/" + "202"
# (7 caracteres mas alla del hallazgo), justo antes del ultimo run "5".
$oldBlock = $d.Range($dateZone.Start, $dateZone.End + 7)
$oldBlock.Delete()
$oldBlock.Collapse(1)
$oldBlock.InsertAfter(" ______________________  ")
$oldBlock.Collapse(0)
$oldBlock.InsertAfter("/202")

# -----------------------------------------------------------------
# 2) Firma con nombre: se agrega " {apellido}" luego de "{nombre}"
#    para que el documento tambien imprima el apellido.
# -----------------------------------------------------------------
$nombre = $d.Content
$foundNombre = $nombre.Find.Execute("{nombre}", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundNombre) { throw "No se encontro el marcador {nombre}." }

$nombre.Collapse(0)
$nombre.Font.Name = "Arial"
$nombre.Font.Bold = $true
$nombre.Font.Size = 10
$nombre.Font.Color = 0
$nombre.InsertAfter(" {apellido}")
